$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (years 2000-2009), shifting rows 12-22 up to become rows 2-12
$ws.Rows("2:11").Delete()

# Add new data rows for 2021 and 2022 at rows 13 and 14
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 74652
$ws.Range("C13").Value = 17072
$ws.Range("D13").Value = 35868
$ws.Range("E13").Value = 21712

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 73351
$ws.Range("C14").Value = 17663
$ws.Range("D14").Value = 34583
$ws.Range("E14").Value = 21105

# Apply the same formatting as the other A column (year) cells to the new rows
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
